$d = $word.ActiveDocument

# Hybrid bold + color (#2C3E50) highlighting for quantitative impact metrics.
# Color value is the decimal (BGR-packed) form of hex 2C3E50 -> RGB(44,62,80)
# as produced by VBA's RGB()/wdColor convention: R + G*256 + B*65536.
$metricColor = 5258796

function Highlight-Terms($paraIndex, $terms) {
    $p = $d.Paragraphs.Item($paraIndex)
    $pStart = $p.Range.Start
    $pEnd = $p.Range.End
    $searchStart = $pStart
    foreach ($term in $terms) {
        $rng = $d.Range($searchStart, $pEnd)
        $found = $rng.Find.Execute($term, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($found) {
            $rng.Font.Bold = 1
            $rng.Font.Color = $metricColor
            $searchStart = $rng.End
        }
    }
}

# --- Professional Experience: Partner - Siege Analytics ---

# "Discovered systematic race coding errors ... from 23% to 64%"
Highlight-Terms 10 @("23%", "64%")

# "Utilized advanced sampling methods ... ±4.2% to ±2.1% ... 71% to 87% ..."
Highlight-Terms 12 @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%", "71%", "87%")

# "Trigonometric algorithm ... reduced mapping costs by 73.5% ... organizations $4.7M ..."
Highlight-Terms 13 @("73.5%", "`$4.7M")

# "Built real-time FEC analysis systems ... valued over $2 trillion"
Highlight-Terms 14 @("`$2")

# --- Professional Experience: Data Products Manager - Helm/Murmuration ---

# "Modernized legacy ETL processes ... reducing processing time by 57%"
Highlight-Terms 20 @("57%")

# --- Key Achievements and Impact ---

# "178% accuracy improvement in racial classification algorithms"
Highlight-Terms 85 @("178%")

# "Algorithmic innovation: ... reducing mapping costs 73.5%"
Highlight-Terms 86 @("73.5%")

# "$4.7M savings enabled nonprofit access"
Highlight-Terms 87 @("`$4.7M")

# "Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"
Highlight-Terms 88 @("12,847")

# "Predictive excellence: ... margin of error from ±4.2% to ±2.1%"
Highlight-Terms 90 @([char]0x00B1 + "4.2%", [char]0x00B1 + "2.1%")

# "Increased voter turnout prediction accuracy from 71% to 87%"
Highlight-Terms 91 @("71%", "87%")
